$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference date format used by the existing "Fecha" column (style s="2")
$dateFormat = $ws.Range("D2").NumberFormat

# Row 2
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C2").Value = 'Ñuble'
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 'Fruta'
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = 'Frutos de pepita'
$ws.Range("I2").Value = 100104003
$ws.Range("J2").Value = 'Membrillo'
$ws.Range("K2").Value = 'Champion'
$ws.Range("D2").Value = 45069
$ws.Range("D2").NumberFormat = $dateFormat
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("Q2").Value = '$/caja 18 kilos empedrada'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 667
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C3").Value = 'Ñuble'
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 'Fruta'
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = 'Frutos de pepita'
$ws.Range("I3").Value = 100104003
$ws.Range("J3").Value = 'Membrillo'
$ws.Range("K3").Value = 'Champion'
$ws.Range("D3").Value = 45069
$ws.Range("D3").NumberFormat = $dateFormat
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = '$/caja 18 kilos empedrada'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 556
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C4").Value = 'Ñuble'
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 'Fruta'
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = 'Frutos de pepita'
$ws.Range("I4").Value = 100104003
$ws.Range("J4").Value = 'Membrillo'
$ws.Range("K4").Value = 'Champion'
$ws.Range("D4").Value = 45044
$ws.Range("D4").NumberFormat = $dateFormat
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 18 kilos empedrada'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 722
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C5").Value = 'Ñuble'
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = 'Frutos de pepita'
$ws.Range("I5").Value = 100104003
$ws.Range("J5").Value = 'Membrillo'
$ws.Range("K5").Value = 'Champion'
$ws.Range("D5").Value = 45044
$ws.Range("D5").NumberFormat = $dateFormat
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/caja 18 kilos empedrada'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C6").Value = 'Ñuble'
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = 'Frutos de pepita'
$ws.Range("I6").Value = 100104003
$ws.Range("J6").Value = 'Membrillo'
$ws.Range("K6").Value = 'Champion'
$ws.Range("D6").Value = 45070
$ws.Range("D6").NumberFormat = $dateFormat
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = '$/caja 18 kilos empedrada'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 556
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C7").Value = 'Ñuble'
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 'Fruta'
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = 'Frutos de pepita'
$ws.Range("I7").Value = 100104003
$ws.Range("J7").Value = 'Membrillo'
$ws.Range("K7").Value = 'Champion'
$ws.Range("D7").Value = 45049
$ws.Range("D7").NumberFormat = $dateFormat
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("Q7").Value = '$/caja 18 kilos empedrada'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 722
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C8").Value = 'Ñuble'
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = 'Frutos de pepita'
$ws.Range("I8").Value = 100104003
$ws.Range("J8").Value = 'Membrillo'
$ws.Range("K8").Value = 'Champion'
$ws.Range("D8").Value = 45049
$ws.Range("D8").NumberFormat = $dateFormat
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = '$/caja 18 kilos empedrada'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 667
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C9").Value = 'Ñuble'
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 'Fruta'
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = 'Frutos de pepita'
$ws.Range("I9").Value = 100104003
$ws.Range("J9").Value = 'Membrillo'
$ws.Range("K9").Value = 'Champion'
$ws.Range("D9").Value = 44699
$ws.Range("D9").NumberFormat = $dateFormat
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/caja 15 kilos granel'
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 867
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C10").Value = 'Ñuble'
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 'Fruta'
$ws.Range("G10").Value = 100104
$ws.Range("H10").Value = 'Frutos de pepita'
$ws.Range("I10").Value = 100104003
$ws.Range("J10").Value = 'Membrillo'
$ws.Range("K10").Value = 'Champion'
$ws.Range("D10").Value = 44699
$ws.Range("D10").NumberFormat = $dateFormat
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 11500
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("R10").Value = 'Provincia de Curicó'
$ws.Range("S10").Value = 767
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C11").Value = 'Ñuble'
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 'Fruta'
$ws.Range("G11").Value = 100104
$ws.Range("H11").Value = 'Frutos de pepita'
$ws.Range("I11").Value = 100104003
$ws.Range("J11").Value = 'Membrillo'
$ws.Range("K11").Value = 'Champion'
$ws.Range("D11").Value = 45033
$ws.Range("D11").NumberFormat = $dateFormat
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 13000
$ws.Range("Q11").Value = '$/caja 18 kilos empedrada'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 722
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C12").Value = 'Ñuble'
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 'Fruta'
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = 'Frutos de pepita'
$ws.Range("I12").Value = 100104003
$ws.Range("J12").Value = 'Membrillo'
$ws.Range("K12").Value = 'Champion'
$ws.Range("D12").Value = 45033
$ws.Range("D12").NumberFormat = $dateFormat
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("Q12").Value = '$/caja 18 kilos empedrada'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 667
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C13").Value = 'Ñuble'
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100104
$ws.Range("H13").Value = 'Frutos de pepita'
$ws.Range("I13").Value = 100104003
$ws.Range("J13").Value = 'Membrillo'
$ws.Range("K13").Value = 'Champion'
$ws.Range("D13").Value = 45062
$ws.Range("D13").NumberFormat = $dateFormat
$ws.Range("L13").Value = 'Especial'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 13000
$ws.Range("Q13").Value = '$/caja 18 kilos empedrada'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 722
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C14").Value = 'Ñuble'
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 'Fruta'
$ws.Range("G14").Value = 100104
$ws.Range("H14").Value = 'Frutos de pepita'
$ws.Range("I14").Value = 100104003
$ws.Range("J14").Value = 'Membrillo'
$ws.Range("K14").Value = 'Champion'
$ws.Range("D14").Value = 45062
$ws.Range("D14").NumberFormat = $dateFormat
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = '$/caja 18 kilos empedrada'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 667
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C15").Value = 'Ñuble'
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100104
$ws.Range("H15").Value = 'Frutos de pepita'
$ws.Range("I15").Value = 100104003
$ws.Range("J15").Value = 'Membrillo'
$ws.Range("K15").Value = 'Champion'
$ws.Range("D15").Value = 45021
$ws.Range("D15").NumberFormat = $dateFormat
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 667
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C16").Value = 'Ñuble'
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = 'Frutos de pepita'
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = 'Membrillo'
$ws.Range("K16").Value = 'Champion'
$ws.Range("D16").Value = 45050
$ws.Range("D16").NumberFormat = $dateFormat
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 13000
$ws.Range("P16").Value = 13000
$ws.Range("Q16").Value = '$/caja 18 kilos empedrada'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 722
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C17").Value = 'Ñuble'
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100104
$ws.Range("H17").Value = 'Frutos de pepita'
$ws.Range("I17").Value = 100104003
$ws.Range("J17").Value = 'Membrillo'
$ws.Range("K17").Value = 'Champion'
$ws.Range("D17").Value = 45050
$ws.Range("D17").NumberFormat = $dateFormat
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = '$/caja 18 kilos empedrada'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 667
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C18").Value = 'Ñuble'
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 'Fruta'
$ws.Range("G18").Value = 100104
$ws.Range("H18").Value = 'Frutos de pepita'
$ws.Range("I18").Value = 100104003
$ws.Range("J18").Value = 'Membrillo'
$ws.Range("K18").Value = 'Champion'
$ws.Range("D18").Value = 45043
$ws.Range("D18").NumberFormat = $dateFormat
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 13000
$ws.Range("Q18").Value = '$/caja 18 kilos empedrada'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 722
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C19").Value = 'Ñuble'
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 'Fruta'
$ws.Range("G19").Value = 100104
$ws.Range("H19").Value = 'Frutos de pepita'
$ws.Range("I19").Value = 100104003
$ws.Range("J19").Value = 'Membrillo'
$ws.Range("K19").Value = 'Champion'
$ws.Range("D19").Value = 45043
$ws.Range("D19").NumberFormat = $dateFormat
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("Q19").Value = '$/caja 18 kilos empedrada'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 667
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C20").Value = 'Ñuble'
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 'Fruta'
$ws.Range("G20").Value = 100104
$ws.Range("H20").Value = 'Frutos de pepita'
$ws.Range("I20").Value = 100104003
$ws.Range("J20").Value = 'Membrillo'
$ws.Range("K20").Value = 'Champion'
$ws.Range("D20").Value = 45020
$ws.Range("D20").NumberFormat = $dateFormat
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 12000
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 667
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C21").Value = 'Ñuble'
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 'Fruta'
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = 'Frutos de pepita'
$ws.Range("I21").Value = 100104003
$ws.Range("J21").Value = 'Membrillo'
$ws.Range("K21").Value = 'Champion'
$ws.Range("D21").Value = 45040
$ws.Range("D21").NumberFormat = $dateFormat
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("Q21").Value = '$/caja 18 kilos empedrada'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 722
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C22").Value = 'Ñuble'
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 'Fruta'
$ws.Range("G22").Value = 100104
$ws.Range("H22").Value = 'Frutos de pepita'
$ws.Range("I22").Value = 100104003
$ws.Range("J22").Value = 'Membrillo'
$ws.Range("K22").Value = 'Champion'
$ws.Range("D22").Value = 45040
$ws.Range("D22").NumberFormat = $dateFormat
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("Q22").Value = '$/caja 18 kilos empedrada'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 667
$ws.Range("T22").Value = 18

